$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (e.g. "579.94", "66.939.98")
# that must remain plain text exactly as scraped. Force text format before
# assigning, then restore default formatting/style so no stray styling is left
# on the cell (matches the source workbook, which keeps these cells unstyled).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "66.939.98"
$ws.Range("E2").Value = "  +4.15%  "
Set-TextValue "D3" "3.268.01"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "579.94"
$ws.Range("E5").Value = "  +3.05%  "
Set-TextValue "D6" "177.05"
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.606"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.03%  "
Set-TextValue "D9" "3.263.36"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("E11").Value = "  +1.33%  "
Set-TextValue "D12" "0.409"
$ws.Range("E12").Value = "  +3.49%  "
Set-TextValue "D13" "3.827.95"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("E14").Value = "  +1.06%  "
Set-TextValue "D15" "28.10"
$ws.Range("E15").Value = "  +2.60%  "
Set-TextValue "D16" "66.934.49"
$ws.Range("E16").Value = "  +4.19%  "
$ws.Range("E17").Value = "  +3.15%  "
Set-TextValue "D18" "3.259.64"
$ws.Range("E18").Value = "  +2.16%  "
Set-TextValue "D19" "5.82"
$ws.Range("E19").Value = "  +2.78%  "
Set-TextValue "D20" "13.41"
$ws.Range("E20").Value = "  +2.42%  "
Set-TextValue "D21" "370.10"
$ws.Range("E21").Value = "  +5.14%  "
Set-TextValue "D22" "7.62"
$ws.Range("E22").Value = "  +6.03%  "
$ws.Range("E23").Value = "  +0.29%  "
Set-TextValue "D24" "70.83"
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("E25").Value = "  +1.00%  "
Set-TextValue "D26" "3.395.55"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("E27").Value = "  +0.26%  "
Set-TextValue "D28" "9.78"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  +2.25%  "
Set-TextValue "D30" "0.999"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +4.67%  "
Set-TextValue "D32" "5.63"
$ws.Range("E32").Value = "  -0.39%  "
Set-TextValue "D33" "22.60"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("E34").Value = "  +0.00%  "
Set-TextValue "D35" "1.25"
$ws.Range("E35").Value = "  +3.73%  "
Set-TextValue "D36" "6.78"
$ws.Range("E36").Value = "  +1.95%  "
Set-TextValue "D37" "170.66"
$ws.Range("E37").Value = "  +9.40%  "
$ws.Range("E38").Value = "  +4.33%  "
Set-TextValue "D39" "0.860"
$ws.Range("E39").Value = "  +6.18%  "
Set-TextValue "D41" "27.17"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D42" "2.756.61"
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.56"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("E44").Value = "  +6.83%  "
$ws.Range("E45").Value = "  +3.93%  "
Set-TextValue "D46" "342.61"
$ws.Range("E46").Value = "  +4.36%  "
Set-TextValue "D47" "40.43"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("E48").Value = "  +3.24%  "
Set-TextValue "D49" "24.71"
$ws.Range("E49").Value = "  +4.52%  "
Set-TextValue "D50" "0.0279"
$ws.Range("E50").Value = "  +2.97%  "
$ws.Range("E51").Value = "  +3.00%  "
